$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so values such as
# "28.20" or "9.10" are not coerced into numbers and lose trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "46.564.99"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.580.82"
$ws.Range("E3").Value = "  +9.63%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "305.69"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "100.12"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("D7").Value = "0.601"
$ws.Range("E7").Value = "  +5.67%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  +12.01%  "
$ws.Range("D10").Value = "38.57"
$ws.Range("E10").Value = "  +11.98%  "
$ws.Range("D11").Value = "0.0839"
$ws.Range("E11").Value = "  +5.18%  "
$ws.Range("D12").Value = "8.16"
$ws.Range("E12").Value = "  +14.61%  "
$ws.Range("D13").Value = "2.971.50"
$ws.Range("E13").Value = "  +9.64%  "
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "2.578.15"
$ws.Range("E15").Value = "  +9.41%  "
$ws.Range("D16").Value = "0.906"
$ws.Range("E16").Value = "  +12.08%  "
$ws.Range("E17").Value = "  +9.56%  "
$ws.Range("D18").Value = "46.678.84"
$ws.Range("D19").Value = "13.41"
$ws.Range("E19").Value = "  +5.69%  "
$ws.Range("E20").Value = "  +4.04%  "
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").Value = "  +10.98%  "
$ws.Range("D22").Value = "70.68"
$ws.Range("E22").Value = "  +5.03%  "
$ws.Range("D23").Value = "253.78"
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  +14.35%  "
$ws.Range("D26").Value = "28.20"
$ws.Range("E26").Value = "  +35.28%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "10.49"
$ws.Range("E28").Value = "  +6.88%  "
$ws.Range("E29").Value = "  +4.39%  "
$ws.Range("D30").Value = "39.62"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("D32").Value = "6.14"
$ws.Range("E32").Value = "  +11.23%  "
$ws.Range("E33").Value = "  +24.36%  "
$ws.Range("D34").Value = "2.92"
$ws.Range("E34").Value = "  +5.40%  "
$ws.Range("E35").Value = "  +7.19%  "
$ws.Range("D36").Value = "150.05"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("E37").Value = "  +4.29%  "
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("D39").Value = "4.20"
$ws.Range("E39").Value = "  +6.31%  "
$ws.Range("D40").Value = "15.73"
$ws.Range("E40").Value = "  +4.26%  "
$ws.Range("D41").Value = "3.62"
$ws.Range("E41").Value = "  +12.51%  "
$ws.Range("E42").Value = "  +7.48%  "
$ws.Range("D43").Value = "2.018.62"
$ws.Range("E43").Value = "  +7.55%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "18.23"
$ws.Range("E45").Value = "  +28.04%  "
$ws.Range("D46").Value = "91.95"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "9.10"
$ws.Range("E48").Value = "  +10.06%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "109.02"
$ws.Range("E49").Value = "  +11.52%  "
$ws.Range("E50").Value = "  +7.66%  "
$ws.Range("D51").Value = "2.830.94"
$ws.Range("E51").Value = "  +9.68%  "
